# "module col dynamic change"
#
# Updates the MonthByDate schedule-notify test data (row 2) on every
# worksheet: new mobile/phone numbers, a new Enquiry/User1 date (and
# date+time), a bumped User1_MessageId counter, and a new User1RecId.
#
# Cells such as phone numbers / dates are stored as literal TEXT in this
# workbook (not numbers/dates), even though their content looks numeric.
# Assigning a numeric- or date-looking string straight to Range.Value
# would make Excel coerce it into a number/date and could even strip
# leading zeros (e.g. "0631225414" -> 631225414). To avoid that we:
#   1. Temporarily force the cell's NumberFormat to "@" (Text) so the
#      value is stored verbatim as text.
#   2. Restore the cell's original look by pasting the formats (only)
#      from a same-style donor cell (column A header style used by all
#      of these row-2 cells), so number formatting / style stay intact.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $ws,
        [string]$cellAddr,
        [string]$text,
        $formatDonor
    )

    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $text

    $formatDonor.Copy()
    $r.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

foreach ($ws in $wb.Worksheets) {
    $donor = $ws.Range("A1")

    Set-TextValue $ws "F2"  "9930402232" $donor
    Set-TextValue $ws "AE2" "0631225414" $donor
    Set-TextValue $ws "AT2" "1514597022" $donor
    Set-TextValue $ws "AX2" "8856092863" $donor

    if ($ws.Name -eq "Sheet1") {
        Set-TextValue $ws "N2"  "2024-01-24" $donor
        Set-TextValue $ws "P2"  "2024-01-24 03:49:03 PM" $donor
        Set-TextValue $ws "AC2" "2024-01-24" $donor
        Set-TextValue $ws "AK2" "2" $donor
        Set-TextValue $ws "AN2" "93950" $donor
    }
}

$excel.CutCopyMode = $false
